$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") "42.752.19"
Set-TextValue $ws.Range("E2") "  -2.47%  "
# Row 3
Set-TextValue $ws.Range("D3") "2.248.07"
Set-TextValue $ws.Range("E3") "  -2.08%  "
# Row 4
Set-TextValue $ws.Range("E4") "  -0.24%  "
# Row 5
Set-TextValue $ws.Range("D5") "116.97"
Set-TextValue $ws.Range("E5") "  +1.42%  "
# Row 6
Set-TextValue $ws.Range("D6") "292.35"
Set-TextValue $ws.Range("E6") "  +9.65%  "
# Row 7
Set-TextValue $ws.Range("D7") "0.634"
Set-TextValue $ws.Range("E7") "  -1.60%  "
# Row 8
Set-TextValue $ws.Range("E8") "  -0.39%  "
# Row 9
Set-TextValue $ws.Range("D9") "0.620"
Set-TextValue $ws.Range("E9") "  +0.68%  "
# Row 10
Set-TextValue $ws.Range("D10") "46.66"
Set-TextValue $ws.Range("E10") "  -1.94%  "
# Row 11
Set-TextValue $ws.Range("E11") "  -0.29%  "
# Row 12
Set-TextValue $ws.Range("D12") "9.18"
Set-TextValue $ws.Range("E12") "  -0.65%  "
# Row 13
Set-TextValue $ws.Range("E13") "  -2.65%  "
# Row 14
Set-TextValue $ws.Range("D14") "15.53"
Set-TextValue $ws.Range("E14") "  +0.36%  "
# Row 15
Set-TextValue $ws.Range("D15") "0.897"
Set-TextValue $ws.Range("E15") "  +2.39%  "
# Row 16
Set-TextValue $ws.Range("D16") "2.587.18"
Set-TextValue $ws.Range("E16") "  -2.08%  "
# Row 17
Set-TextValue $ws.Range("D17") "2.242.47"
Set-TextValue $ws.Range("E17") "  -2.42%  "
# Row 18
Set-TextValue $ws.Range("D18") "42.784.66"
Set-TextValue $ws.Range("E18") "  -2.20%  "
# Row 19
Set-TextValue $ws.Range("D19") "7.65"
Set-TextValue $ws.Range("E19") "  +12.61%  "
# Row 20
Set-TextValue $ws.Range("E20") "  -1.57%  "
# Row 21
Set-TextValue $ws.Range("E21") "  +1.72%  "
# Row 22
Set-TextValue $ws.Range("D22") "3.45"
Set-TextValue $ws.Range("E22") "  +18.45%  "
# Row 23
Set-TextValue $ws.Range("D23") "2.38"
Set-TextValue $ws.Range("E23") "  -2.52%  "
# Row 24
Set-TextValue $ws.Range("D24") "233.66"
Set-TextValue $ws.Range("E24") "  -1.37%  "
# Row 25
Set-TextValue $ws.Range("D25") "9.38"
Set-TextValue $ws.Range("E25") "  -1.22%  "
# Row 26
Set-TextValue $ws.Range("D26") "12.24"
Set-TextValue $ws.Range("E26") "  +5.18%  "
# Row 27
Set-TextValue $ws.Range("D27") "0.999"
Set-TextValue $ws.Range("E27") "  -1.78%  "
# Row 28
Set-TextValue $ws.Range("D28") "40.47"
Set-TextValue $ws.Range("E28") "  -2.92%  "
# Row 29
Set-TextValue $ws.Range("B29") "WEMIXToken"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D29") "3.28"
Set-TextValue $ws.Range("E29") "  -3.22%  "
# Row 30
Set-TextValue $ws.Range("B30") "Toncoin"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D30") "2.19"
Set-TextValue $ws.Range("E30") "  -2.33%  "
# Row 31
Set-TextValue $ws.Range("D31") "175.67"
# Row 32
Set-TextValue $ws.Range("D32") "21.35"
Set-TextValue $ws.Range("E32") "  -2.18%  "
# Row 33
Set-TextValue $ws.Range("D33") "0.0915"
Set-TextValue $ws.Range("E33") "  +0.45%  "
# Row 34
Set-TextValue $ws.Range("D34") "4.59"
Set-TextValue $ws.Range("E34") "  +17.70%  "
# Row 35
Set-TextValue $ws.Range("E35") "  -0.29%  "
# Row 36
Set-TextValue $ws.Range("E36") "  -1.68%  "
# Row 37
Set-TextValue $ws.Range("E37") "  +1.05%  "
# Row 38
Set-TextValue $ws.Range("E38") "  -1.58%  "
# Row 39
Set-TextValue $ws.Range("E39") "  +0.35%  "
# Row 40
Set-TextValue $ws.Range("D40") "2.63"
# Row 41
Set-TextValue $ws.Range("D41") "72.92"
Set-TextValue $ws.Range("E41") "  -2.83%  "
# Row 42
Set-TextValue $ws.Range("E42") "  +1.09%  "
# Row 43
Set-TextValue $ws.Range("D43") "13.60"
Set-TextValue $ws.Range("E43") "  -6.32%  "
# Row 44
Set-TextValue $ws.Range("E44") "  -0.20%  "
# Row 45
Set-TextValue $ws.Range("E45") "  -1.84%  "
# Row 46
Set-TextValue $ws.Range("D46") "5.61"
Set-TextValue $ws.Range("E46") "  -7.55%  "
# Row 47
Set-TextValue $ws.Range("D47") "1.32"
Set-TextValue $ws.Range("E47") "  +2.23%  "
# Row 48
Set-TextValue $ws.Range("D48") "105.92"
Set-TextValue $ws.Range("E48") "  +5.44%  "
# Row 49
Set-TextValue $ws.Range("D49") "8.65"
Set-TextValue $ws.Range("E49") "  +0.22%  "
# Row 50
Set-TextValue $ws.Range("B50") "Cronos"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D50") "0.0992"
Set-TextValue $ws.Range("E50") "  -1.05%  "
# Row 51
Set-TextValue $ws.Range("D51") "0.472"
Set-TextValue $ws.Range("E51") "  +7.10%  "
